# edit.ps1 - apply the Test_PGS_getBadges.docx changes via Word COM interop
#
# Summary of changes:
#  1) "Plan de test - Fonction [getBadges]" : merge the two runs that spell
#     out "getBadges" (wrapped in a spellStart/spellEnd proofErr pair) into
#     a single run, dropping the proofErr markers.
#  2) " : T7.1" -> " : T7"
#  3) "routes/pgs.py - Fonction [getBadges]" : same run/proofErr merge as (1).
#  4) "actif ([boolean])" : same merge, dropping proofErr around "boolean".
#  5) "[creation] (date)" : same merge, dropping proofErr around "creation".
#  6) "Colonnes uid, actif, [creation], id_utilisateur ..." : same merge.
#  7) word/styles.xml: the "Policepardfaut" (Default Paragraph Font) style
#     gains a <w:semiHidden/> flag.

$d = $word.ActiveDocument

function Get-ParagraphIndexContaining([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# (1) Title paragraph: "Plan de test - Fonction getBadges"
#     proofErr(spellStart) precedes "getBadges" (a run exists before it,
#     so only a trailing dummy char is needed to give the replace range a
#     run on both sides of the proofErr pair).
# ---------------------------------------------------------------------
$idx = Get-ParagraphIndexContaining("Plan de test")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$collapsed = $d.Range($r.End - 1, $r.End - 1)
$collapsed.InsertAfter("ZZMARKER")
$d.Content.Find.Execute("Plan de test – Fonction getBadgesZZMARKER", $true, $false, $false, $false, $false, $true, 1, $false, "Plan de test – Fonction getBadges", 2) | Out-Null

# ---------------------------------------------------------------------
# (2) "Numero : T7.1" -> "Numero : T7"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" : T7.1", $true, $false, $false, $false, $false, $true, 1, $false, " : T7", 2) | Out-Null

# ---------------------------------------------------------------------
# (3) "routes/pgs.py - Fonction getBadges" (proofErr spellEnd is the very
#     last element of the paragraph -> needs a trailing dummy char).
# ---------------------------------------------------------------------
$idx = Get-ParagraphIndexContaining("routes/pgs.py")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$collapsed = $d.Range($r.End - 1, $r.End - 1)
$collapsed.InsertAfter("ZZMARKER")
$d.Content.Find.Execute("routes/pgs.py – Fonction getBadgesZZMARKER", $true, $false, $false, $false, $false, $true, 1, $false, "routes/pgs.py – Fonction getBadges", 2) | Out-Null

# ---------------------------------------------------------------------
# (4) "actif (boolean)" - proofErr pair already has runs on both sides.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("actif (boolean)", $true, $false, $false, $false, $false, $true, 1, $false, "actif (boolean)", 2) | Out-Null

# ---------------------------------------------------------------------
# (5) "creation (date)" - proofErr spellStart is the very first element of
#     the paragraph -> needs a leading dummy char.
# ---------------------------------------------------------------------
$idx = Get-ParagraphIndexContaining("(date)")
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$collapsed = $d.Range($r.Start, $r.Start)
$collapsed.InsertBefore("ZZMARKER")
$d.Content.Find.Execute("ZZMARKERcreation (date)", $true, $false, $false, $false, $false, $true, 1, $false, "creation (date)", 2) | Out-Null

# ---------------------------------------------------------------------
# (6) "Colonnes uid, actif, creation, id_utilisateur correctement
#     remplies" - proofErr pair already interior to the paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Colonnes uid, actif, creation, id_utilisateur correctement remplies", $true, $false, $false, $false, $false, $true, 1, $false, "Colonnes uid, actif, creation, id_utilisateur correctement remplies", 2) | Out-Null

# ---------------------------------------------------------------------
# (7) styles.xml: Policepardfaut (Default Paragraph Font) gains
#     <w:semiHidden/>.
# ---------------------------------------------------------------------
$style = $d.Styles("Policepardfaut")
$style.NoProofing = $style.NoProofing
